# Add team record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 55

# --- Header row (row 1): new column headers, matching the style of the
#     existing header cells (bold, centered, thin border). We copy the
#     formatting from the last existing header cell (AC1) onto the new
#     header cells before setting their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2..55): every player row gets the team's season record.
$rowCount = $lastRow - 1
$data = New-Object 'object[,]' $rowCount,3
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i,0] = 89
    $data[$i,1] = 73
    $data[$i,2] = 0
}
$ws.Range("AD2:AF$lastRow").Value = $data
